# Ran DDE model for all insect species
#
# The "Aphis citricola Japan Chiba" entry had not actually been run through
# the DDE model, so that whole row is removed from the results sheet and all
# rows below it shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model results")

# Delete the entire row containing "Aphis citricola Japan Chiba" (row 17).
# This shifts rows 18:31 up to 17:30, so the data for every other species is
# preserved but the sheet now ends at row 30 instead of row 31.
$ws.Rows("17").Delete()

# Leave the selection where Excel would naturally land after deleting the
# row and scrolling to the (now empty) row following the last data row.
$ws.Range("A31").Select()
